$wb = $excel.ActiveWorkbook

# --- Timers sheet: add a new "TIM2 / Task scheduling (32 bit)" row above the
#     existing TIM3 row (new row 5, pushing the rest down by one). ---
$timers = $wb.Worksheets.Item("Timers")
$timers.Rows.Item(5).Insert() | Out-Null
$timers.Range("A5").Value = "TIM2"
$timers.Range("B5").Value = "Task scheduling (32 bit)"

# --- Interrupt Priorities sheet: add a new "Priority 3" section at the
#     bottom with a "TIM2 / Task scheduler timer" entry. ---
$prio = $wb.Worksheets.Item("Interrupt Priorities")
$prio.Range("A21").Value = "Priority 3"
$prio.Range("A21").Style = "Heading 1"
$prio.Range("A23").Value = "TIM2"
$prio.Range("B23").Value = "Task scheduler timer"

# --- Selection / active-sheet bookkeeping to match the saved UI state ---
$timers.Range("B20").Select() | Out-Null

$prio.Activate() | Out-Null
$prio.Range("A24").Select() | Out-Null
